# The deck's design ("Integral") had its 12 theme colours swapped for the
# colours of the stock "Office Theme" palette (the palette already used
# elsewhere in the package, e.g. by the notes master). Font scheme and
# format scheme (fills/lines/effects) are identical between the two
# palettes, so only the colour scheme itself needs to change.
#
# RGB() packs as 0x00BBGGRR, matching PowerPoint's ColorScheme.Colors(n).RGB:
#   1 dk1     = 000000
#   2 lt1     = FFFFFF
#   3 dk2     = 44546A
#   4 lt2     = E7E6E6
#   5 accent1 = 5B9BD5
#   6 accent2 = ED7D31
#   7 accent3 = A5A5A5
#   8 accent4 = FFC000
#   9 accent5 = 4472C4
#  10 accent6 = 70AD47
#  11 hlink   = 0563C1
#  12 folHlink= 954F72

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388     # dk2      44546A
$cs.Colors(4).RGB  = 15132391    # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939    # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501     # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845    # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407       # accent4  FFC000
$cs.Colors(9).RGB  = 12874308    # accent5  4472C4
$cs.Colors(10).RGB = 4697456     # accent6  70AD47
$cs.Colors(11).RGB = 12673797    # hlink    0563C1
$cs.Colors(12).RGB = 7491477     # folHlink 954F72
